# Added placeholder PPCF effectiveness & affected fraction for wasting
# Target sheet: "Interventions for children" (rows 6-11 get new data,
# mirroring the existing two-row-block / three-row-block pattern already
# present in rows 2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions for children")

# ---------------------------------------------------------------------
# 1. Apply number formatting (style) to the new cells by cloning the
#    formats that already exist on the analogous "highlighted" block in
#    row 5 (style ids 62 / 63), before writing any values. This avoids
#    creating any brand-new style table entries.
# ---------------------------------------------------------------------

# Row 6 (Wasting (high) / Affected fraction) - D:F plain highlighted (62),
# G:H right-aligned highlighted (63), matching row 5's pattern.
$ws.Range("D5:F5").Copy()
$ws.Range("D6:F6").PasteSpecial(-4122)
$ws.Range("G5:H5").Copy()
$ws.Range("G6:H6").PasteSpecial(-4122)

# Row 7 (Effectiveness mortality) - D:H all style 62.
$ws.Range("D5").Copy()
$ws.Range("D7:H7").PasteSpecial(-4122)

# Row 8 (Effectiveness incidence) - C gets style 4 (like C4), D:H style 62.
$ws.Range("C4").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D8:H8").PasteSpecial(-4122)

# Row 9 (Wasting (moderate) / Affected fraction) - C gets style 4, D:F
# style 62, G:H style 63.
$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D5:F5").Copy()
$ws.Range("D9:F9").PasteSpecial(-4122)
$ws.Range("G5:H5").Copy()
$ws.Range("G9:H9").PasteSpecial(-4122)

# Row 10 (Effectiveness mortality) - D:H all style 62.
$ws.Range("D5").Copy()
$ws.Range("D10:H10").PasteSpecial(-4122)

# Row 11 (Effectiveness incidence) - C gets style 4, D:H style 62.
$ws.Range("C4").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D11:H11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Fill in the cell values/text (reusing existing shared strings).
# ---------------------------------------------------------------------

# Row 6
$ws.Range("A6").Value = "Public provision of complementary foods"
$ws.Range("B6").Value = "Wasting (high)"
$ws.Range("C6").Value = "Affected fraction"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.335
$ws.Range("G6").Value = 0.335
$ws.Range("H6").Value = 0.335

# Row 7
$ws.Range("C7").Value = "Effectiveness mortality"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.3
$ws.Range("G7").Value = 0.3
$ws.Range("H7").Value = 0.3

# Row 8
$ws.Range("C8").Value = "Effectiveness incidence"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.62
$ws.Range("G8").Value = 0.62
$ws.Range("H8").Value = 0.62

# Row 9
$ws.Range("B9").Value = "Wasting (moderate)"
$ws.Range("C9").Value = "Affected fraction"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0.335
$ws.Range("G9").Value = 0.335
$ws.Range("H9").Value = 0.335

# Row 10
$ws.Range("C10").Value = "Effectiveness mortality"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.3
$ws.Range("G10").Value = 0.3
$ws.Range("H10").Value = 0.3

# Row 11
$ws.Range("C11").Value = "Effectiveness incidence"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.62
$ws.Range("G11").Value = 0.62
$ws.Range("H11").Value = 0.62

# ---------------------------------------------------------------------
# 3. Restore view-state: selection on each touched sheet, with
#    "Interventions for children" left as the active tab/selection.
# ---------------------------------------------------------------------

$wsBaseline = $wb.Worksheets.Item("Baseline year demographics")
$wsBaseline.Range("D10").Select() | Out-Null

$wsCost = $wb.Worksheets.Item("Interventions cost and coverage")
$wsCost.Range("A6").Select() | Out-Null

$ws.Activate()
$ws.Range("R51").Select() | Out-Null
